$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '27.292.14'
Set-TextValue $ws.Range('E2') '  +2.23%  '

Set-TextValue $ws.Range('D3') '1.820.02'
Set-TextValue $ws.Range('E3') '  +1.33%  '

Set-TextValue $ws.Range('D4') '1.000'
Set-TextValue $ws.Range('E4') '  -0.04%  '

Set-TextValue $ws.Range('D5') '313.57'
Set-TextValue $ws.Range('E5') '  +1.41%  '

Set-TextValue $ws.Range('E6') '  -0.04%  '

Set-TextValue $ws.Range('D7') '0.4642'
Set-TextValue $ws.Range('E7') '  +4.70%  '

Set-TextValue $ws.Range('D8') '0.3767'
Set-TextValue $ws.Range('E8') '  +2.35%  '

Set-TextValue $ws.Range('D9') '0.07430'
Set-TextValue $ws.Range('E9') '  +1.03%  '

Set-TextValue $ws.Range('D10') '0.8709'
Set-TextValue $ws.Range('E10') '  +1.25%  '

Set-TextValue $ws.Range('E11') '  +0.02%  '

Set-TextValue $ws.Range('D12') '1.822.34'
Set-TextValue $ws.Range('E12') '  +1.42%  '

Set-TextValue $ws.Range('D13') '6.677'
Set-TextValue $ws.Range('E13') '  +0.86%  '

Set-TextValue $ws.Range('E14') '  +2.59%  '

Set-TextValue $ws.Range('D15') '0.07098'
Set-TextValue $ws.Range('E15') '  +0.45%  '

Set-TextValue $ws.Range('D16') '92.09'
Set-TextValue $ws.Range('E16') '  +0.33%  '

Set-TextValue $ws.Range('D17') '1.001'
Set-TextValue $ws.Range('E17') '  -0.02%  '

Set-TextValue $ws.Range('D18') '0.000008759'
Set-TextValue $ws.Range('E18') '  +0.85%  '

Set-TextValue $ws.Range('E19') '  -0.06%  '

Set-TextValue $ws.Range('E20') '  +0.90%  '

Set-TextValue $ws.Range('D21') '27.295.74'
Set-TextValue $ws.Range('E21') '  +2.17%  '

Set-TextValue $ws.Range('D22') '5.308'
Set-TextValue $ws.Range('E22') '  +2.83%  '

Set-TextValue $ws.Range('D23') '10.94'
Set-TextValue $ws.Range('E23') '  +1.12%  '

Set-TextValue $ws.Range('D24') '2.050.86'
Set-TextValue $ws.Range('E24') '  +1.51%  '

Set-TextValue $ws.Range('D25') '1.937'
Set-TextValue $ws.Range('E25') '  -2.07%  '

Set-TextValue $ws.Range('D26') '151.73'
Set-TextValue $ws.Range('E26') '  -0.20%  '

Set-TextValue $ws.Range('D27') '2.243'
Set-TextValue $ws.Range('E27') '  +2.93%  '

Set-TextValue $ws.Range('D28') '18.61'
Set-TextValue $ws.Range('E28') '  +1.01%  '

Set-TextValue $ws.Range('D29') '5.279'
Set-TextValue $ws.Range('E29') '  +1.90%  '

Set-TextValue $ws.Range('D30') '117.25'
Set-TextValue $ws.Range('E30') '  +0.08%  '

Set-TextValue $ws.Range('D31') '0.08913'
Set-TextValue $ws.Range('E31') '  +1.55%  '

Set-TextValue $ws.Range('D32') '0.7818'
Set-TextValue $ws.Range('E32') '  +5.70%  '

Set-TextValue $ws.Range('E33') '  +2.57%  '

Set-TextValue $ws.Range('D34') '4.519'
Set-TextValue $ws.Range('E34') '  +1.70%  '

Set-TextValue $ws.Range('D35') '2.923'
Set-TextValue $ws.Range('E35') '  +0.72%  '

Set-TextValue $ws.Range('D36') '1.000'
Set-TextValue $ws.Range('E36') '  +0.00%  '

Set-TextValue $ws.Range('D37') '1.101'
Set-TextValue $ws.Range('E37') '  +1.60%  '

Set-TextValue $ws.Range('D38') '0.01970'
Set-TextValue $ws.Range('E38') '  +0.67%  '

Set-TextValue $ws.Range('D39') '0.05260'
Set-TextValue $ws.Range('E39') '  +1.33%  '

Set-TextValue $ws.Range('D40') '7.303'
Set-TextValue $ws.Range('E40') '  +4.79%  '

$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range('D41') '0.5291'
Set-TextValue $ws.Range('E41') '  +0.78%  '

$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range('D42') '2.892'
Set-TextValue $ws.Range('E42') '  +2.45%  '

$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D43') '2.358'
Set-TextValue $ws.Range('E43') '  +19.30%  '

Set-TextValue $ws.Range('D44') '0.1689'
Set-TextValue $ws.Range('E44') '  +0.38%  '

Set-TextValue $ws.Range('D45') '8.602'
Set-TextValue $ws.Range('E45') '  +1.83%  '

Set-TextValue $ws.Range('D46') '0.5039'
Set-TextValue $ws.Range('E46') '  -0.33%  '

Set-TextValue $ws.Range('D47') '10.51'
Set-TextValue $ws.Range('E47') '  +0.68%  '

Set-TextValue $ws.Range('D48') '105.53'
Set-TextValue $ws.Range('E48') '  +0.44%  '

Set-TextValue $ws.Range('D49') '1.674'
Set-TextValue $ws.Range('E49') '  +0.53%  '

Set-TextValue $ws.Range('E50') '  +0.01%  '

Set-TextValue $ws.Range('D51') '0.06324'
Set-TextValue $ws.Range('E51') '  +0.57%  '
